$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert 5 new rows AFTER row 40 (row 40 itself keeps its own data: MARIA ALEJANDRA / 2508) ---
$ws.Range("A41:A45").EntireRow.Insert()

# New rows 41-45 (and row 40, no longer the last table row) get "normal data row" formatting,
# matching the look of rows 16-39 - copy the format from row 39.
$ws.Range("B39:J39").Copy()
$ws.Range("B40:J45").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# --- 2. Fill the new rows with data for period 2509 ---
# Row 41: MELVA CASTRO TORRES
$ws.Cells.Item(41,2).Value = "CC"
$ws.Cells.Item(41,3).Value = "32908729"
$ws.Cells.Item(41,4).Value = "MELVA CASTRO TORRES"
$ws.Cells.Item(41,5).Value = "2509"
$ws.Cells.Item(41,6).Value = 56940
$ws.Cells.Item(41,7).Value = 1423500

# Row 42: DIOSELINA ROMERO AUDIVET
$ws.Cells.Item(42,2).Value = "CC"
$ws.Cells.Item(42,3).Value = "45483943"
$ws.Cells.Item(42,4).Value = "DIOSELINA ROMERO AUDIVET"
$ws.Cells.Item(42,5).Value = "2509"
$ws.Cells.Item(42,6).Value = 56940
$ws.Cells.Item(42,7).Value = 1423500

# Row 43: MANUEL FEDERICO ROMERO VARGAS
$ws.Cells.Item(43,2).Value = "CC"
$ws.Cells.Item(43,3).Value = "73127790"
$ws.Cells.Item(43,4).Value = "MANUEL FEDERICO ROMERO VARGAS"
$ws.Cells.Item(43,5).Value = "2509"
$ws.Cells.Item(43,6).Value = 56940
$ws.Cells.Item(43,7).Value = 1423500

# Row 44: EMILSE CASSIANIS MIRANDA (new worker)
$ws.Cells.Item(44,2).Value = "CC"
$ws.Cells.Item(44,3).Value = "22801691"
$ws.Cells.Item(44,4).Value = "EMILSE CASSIANIS MIRANDA"
$ws.Cells.Item(44,5).Value = "2509"
$ws.Cells.Item(44,6).Value = 37960
$ws.Cells.Item(44,7).Value = 1423500

# Row 45: MARIA ALEJANDRA GAMARRA ALVAREZ (last row of the table)
$ws.Cells.Item(45,2).Value = "CC"
$ws.Cells.Item(45,3).Value = "1143362154"
$ws.Cells.Item(45,4).Value = "MARIA ALEJANDRA GAMARRA ALVAREZ"
$ws.Cells.Item(45,5).Value = "2509"
$ws.Cells.Item(45,6).Value = 56940
$ws.Cells.Item(45,7).Value = 1423500

# --- 3. Center-align the "Periodo Mora" column (E) for all data rows ---
$ws.Range("E16:E45").HorizontalAlignment = -4108 # xlCenter

# --- 4. Give the new last row (45) the distinct "bottom of table" border treatment ---
$ws.Range("B45:J45").Borders.Item(9).Color = 0   # xlEdgeBottom, explicit black
$ws.Range("B45:J45").Borders.Item(9).Weight = 2  # xlThin

# --- 5. Update summary counters ---
$ws.Cells.Item(13,3).Value = 8   # Cant. Trabajadores
$ws.Cells.Item(13,6).Value = 11  # Cant. Periodos
$ws.Cells.Item(11,5).Value = 1566552 # VALOR MORA total

Write-Host $ws.UsedRange.Address()
